# Apply changes: Added ExitOnNextEntry mode and fixed bugs on trade entries for ScalpEmaRsiAdx.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("ListOfValues")

# --- Update existing row 2 (Test # 1) ---
$ws1.Range("D2").Value = 44550
$ws1.Range("H2").Value = 4
$ws1.Range("I2").Value = 4
$ws1.Range("J2").Value = "ScalpEmaRsiAdx"
$ws1.Range("K2").Value = "ExitOnNextEntry"

# --- Add new rows 3-6, cloning the formatting from row 2 first ---
$ws1.Range("A2:K2").Copy()
$ws1.Range("A3:K3").PasteSpecial(-4122)
$ws1.Range("A4:K4").PasteSpecial(-4122)
$ws1.Range("A5:K5").PasteSpecial(-4122)
$ws1.Range("A6:K6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 3: Test # 2
$ws1.Range("A3").Value = 2
$ws1.Range("B3").Value = "Bybit"
$ws1.Range("C3").Value = "BTCUSDT"
$ws1.Range("D3").Value = 44550
$ws1.Range("E3").Value = 44593
$ws1.Range("F3").Value = "3m"
$ws1.Range("G3").Value = 10000
$ws1.Range("H3").Value = 3
$ws1.Range("I3").Value = 3
$ws1.Range("J3").Value = "ScalpEmaRsiAdx"
$ws1.Range("K3").Value = "ExitOnNextEntry"

# Row 4: Test # 3
$ws1.Range("A4").Value = 3
$ws1.Range("B4").Value = "Bybit"
$ws1.Range("C4").Value = "BTCUSDT"
$ws1.Range("D4").Value = 44550
$ws1.Range("E4").Value = 44593
$ws1.Range("F4").Value = "3m"
$ws1.Range("G4").Value = 10000
$ws1.Range("H4").Value = 2
$ws1.Range("I4").Value = 2
$ws1.Range("J4").Value = "ScalpEmaRsiAdx"
$ws1.Range("K4").Value = "ExitOnNextEntry"

# Row 5: Test # 4
$ws1.Range("A5").Value = 4
$ws1.Range("B5").Value = "Bybit"
$ws1.Range("C5").Value = "BTCUSDT"
$ws1.Range("D5").Value = 44550
$ws1.Range("E5").Value = 44593
$ws1.Range("F5").Value = "3m"
$ws1.Range("G5").Value = 10000
$ws1.Range("H5").Value = 1
$ws1.Range("I5").Value = 1
$ws1.Range("J5").Value = "ScalpEmaRsiAdx"
$ws1.Range("K5").Value = "ExitOnNextEntry"

# Row 6: Test # 5 (replaces the old blank placeholder row 6)
$ws1.Range("A6").Value = 5
$ws1.Range("B6").Value = "Bybit"
$ws1.Range("C6").Value = "BTCUSDT"
$ws1.Range("D6").Value = 44550
$ws1.Range("E6").Value = 44593
$ws1.Range("F6").Value = "3m"
$ws1.Range("G6").Value = 10000
$ws1.Range("H6").Value = 1
$ws1.Range("I6").Value = 1
$ws1.Range("J6").Value = "ScalpEmaRsiAdx"
$ws1.Range("K6").Value = "FixedPCT"

# --- Selections (active cell) ---
# Set ListOfValues' selection first, then return to Sheet1 so Sheet1
# remains the active/selected tab, matching the original tabSelected="1".
$null = $ws2.Range("C3").Select()
$null = $ws1.Select()
$null = $ws1.Range("J10").Select()
